$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: paste a brand-new run of plain text, carrying the same run
# formatting (sz=24 / szCs=24, the "body text" size used throughout this
# document) as the rest of the paragraph. A plain Range.InsertAfter() call
# creates a run with NO rPr at all, and Font.Size only ever stamps <w:sz/>
# (never <w:szCs/>) in this host, so instead we stage the text in a scratch
# spot, clone an existing sz/szCs run onto it via Copy/Paste (which keeps
# full formatting), fix up the staged text, then Copy/Paste *that* into the
# real insertion point. Paste() never silently merges into a neighbouring
# run, so every call below yields its own separate <w:r>, matching how the
# target document was produced.
# ---------------------------------------------------------------------------
function Insert-BodyRun($doc, $pos, $text) {
    $donor = $doc.Range(0, 0)
    $donor.Find.Execute("god ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $donor.Copy()

    $scratchStart = $doc.Paragraphs.Item($doc.Paragraphs.Count).Range.Start
    $scratchPoint = $doc.Range($scratchStart, $scratchStart)
    $scratchPoint.Paste()

    $scratchRange = $doc.Range($scratchStart, $scratchStart + 4)
    $scratchRange.Text = $text
    $scratchRange = $doc.Range($scratchStart, $scratchStart + $text.Length)
    $scratchRange.Copy()
    $scratchRange.Text = ""

    $tgt = $doc.Range($pos, $pos)
    $tgt.Paste()

    return $pos + $text.Length
}

# ---------------------------------------------------------------------------
# Part 1: paragraph "Then I remember that i was near him ( god but don't
# show him ) my hands yea my hands they where bloody then I don't remember
# it exactly I promised to him" gets rewritten/expanded from its last two
# runs onward, and the paragraph is split in two, with a further new
# paragraph appended ("Yes doctor every time when I get lost or feel ").
# ---------------------------------------------------------------------------
$old1 = ") my hands yea my hands they where bloody then I don" + [char]0x2019 + "t remember it exactly I promised to him"
$target = $d.Content
$target.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pos = $target.Start
$target.Text = ""

$pos = Insert-BodyRun $d $pos ") "
$pos = Insert-BodyRun $d $pos "Then I saw "
$pos = Insert-BodyRun $d $pos "my hands "
$pos = Insert-BodyRun $d $pos ", "
$pos = Insert-BodyRun $d $pos "yea my hands they where "
$pos = Insert-BodyRun $d $pos "covered in red blood "
$pos = Insert-BodyRun $d $pos "."
$pos = Insert-BodyRun $d $pos " "
$pos = Insert-BodyRun $d $pos "then "
$pos = Insert-BodyRun $d $pos "I all I remember is"
$pos = Insert-BodyRun $d $pos " "
$pos = Insert-BodyRun $d $pos "my hands where shivering then I hold his (god) hands and made that"
$pos = Insert-BodyRun $d $pos " "
$pos = Insert-BodyRun $d $pos "promise"

# Split into a new paragraph right here.
$breakPoint = $d.Range($pos, $pos)
$breakPoint.InsertParagraphAfter()
$pos = $pos + 1

$pos = Insert-BodyRun $d $pos "Yes "
$pos = Insert-BodyRun $d $pos "doctor every time when I get lost or feel "

# ---------------------------------------------------------------------------
# Part 2: a new paragraph "Pc update " is appended right after "So this was
# ur dream - physiatrists to me ", before the first originally-blank
# paragraph.
# ---------------------------------------------------------------------------
$old2 = "So this was ur dream " + [char]0x2013 + " physiatrists to me "
$target2 = $d.Content
$target2.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endPos = $target2.End

$breakPoint2 = $d.Range($endPos, $endPos)
$breakPoint2.InsertParagraphAfter()
$pos2 = $endPos + 1

$pos2 = Insert-BodyRun $d $pos2 "Pc update "

Write-Output "done"
